$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column from 45181 -> 45182 for all existing data rows (2..132)
for ($r = 2; $r -le 132; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45182
}

# Ensure row 132 has an explicit row height (matches surrounding rows)
$ws.Rows.Item(132).RowHeight = 15

# Add the new row 133 with data for case "A 42742-2023"
$row = 133
$ws.Cells.Item($row, 1).Value2 = "A 42742-2023"
$ws.Cells.Item($row, 2).Value2 = 45181
$ws.Cells.Item($row, 3).Value2 = 45182
$ws.Cells.Item($row, 4).Value2 = "STOCKHOLMS LÄN"
$ws.Cells.Item($row, 5).Value2 = "VALLENTUNA"
$ws.Cells.Item($row, 7).Value2 = 10.4
$ws.Cells.Item($row, 8).Value2 = 0
$ws.Cells.Item($row, 9).Value2 = 0
$ws.Cells.Item($row, 10).Value2 = 0
$ws.Cells.Item($row, 11).Value2 = 0
$ws.Cells.Item($row, 12).Value2 = 0
$ws.Cells.Item($row, 13).Value2 = 0
$ws.Cells.Item($row, 14).Value2 = 0
$ws.Cells.Item($row, 15).Value2 = 0
$ws.Cells.Item($row, 16).Value2 = 0
$ws.Cells.Item($row, 17).Value2 = 0

# Apply same formatting (date format) used by the rest of column B and C
$ws.Cells.Item($row, 2).NumberFormat = $ws.Cells.Item($row - 1, 2).NumberFormat
$ws.Cells.Item($row, 3).NumberFormat = $ws.Cells.Item($row - 1, 3).NumberFormat

# Column R (Artnamn) uses a wrap-text style even when empty
$ws.Cells.Item($row, 18).WrapText = $ws.Cells.Item($row - 1, 18).WrapText

$ws.Rows.Item($row).RowHeight = 15
